# Apply "corregi algoritmo calculo de reactivities" edit.
$wb = $excel.ActiveWorkbook

# --- Sheet3: new reactivity / retention summary table (rows 9-59) ---
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Range("D9").Value = "BA:"
$ws.Range("D10").Value = "Deposition"
$ws.Range("E10").Value = "mass"
$ws.Range("F10").Value = "(g/cm2/year):"
$ws.Range("G10").Value = 12.52
$ws.Range("D11").Value = "Sedimentation"
$ws.Range("E11").Value = "rate"
$ws.Range("F11").Value = "(cm/year):"
$ws.Range("G11").Value = 4.73
$ws.Range("D12").Value = "N:"
$ws.Range("D13").Value = "Deposition"
$ws.Range("E13").Value = "mass"
$ws.Range("F13").Value = "(g/cm2/year):"
$ws.Range("G13").Value = 0.55
$ws.Range("D14").Value = "Sedimentation"
$ws.Range("E14").Value = "rate"
$ws.Range("F14").Value = "(cm/year):"
$ws.Range("G14").Value = 0.21
$ws.Range("D15").Value = "Retention"
$ws.Range("E15").Value = "in"
$ws.Range("F15").Value = "sediment"
$ws.Range("G15").Value = "(%)"
$ws.Range("D16").Value = "BA"
$ws.Range("D17").Value = "Coprostanol"
$ws.Range("E17").Value = 0.0671861287979
$ws.Range("F17").Value = 6.5043660888
$ws.Range("D18").Value = "Epicoprostanol"
$ws.Range("E18").Value = 0.00308612389878
$ws.Range("F18").Value = 40.3871076676
$ws.Range("D19").Value = "Ethylcoprostanol"
$ws.Range("E19").Value = 0.00876993233868
$ws.Range("F19").Value = 7.63645821842
$ws.Range("D20").Value = "Coprostanone"
$ws.Range("E20").Value = 0.00502373538883
$ws.Range("F20").Value = 10.0592682916
$ws.Range("D21").Value = "Coprostane"
$ws.Range("E21").Value = 0.0000113736008665
$ws.Range("F21").Value = 166.334396002
$ws.Range("D22").Value = "b-Sitosterol"
$ws.Range("E22").Value = 0.00384525844516
$ws.Range("F22").Value = 9.75985204793
$ws.Range("D23").Value = "γ-Sitosterol"
$ws.Range("E23").Value = 0.000217173355268
$ws.Range("F23").Value = 0.824022071787
$ws.Range("D24").Value = "Stigmasterol"
$ws.Range("E24").Value = 0.000389117984313
$ws.Range("F24").Value = 10.4198917254
$ws.Range("D25").Value = "Stigmastanol"
$ws.Range("E25").Value = 0.00119501642915
$ws.Range("F25").Value = 14.3889804129
$ws.Range("D26").Value = "Campesterol"
$ws.Range("E26").Value = 0.000400093182583
$ws.Range("F26").Value = 11.5363815572
$ws.Range("D27").Value = "Campestanol"
$ws.Range("E27").Value = 0.00000639074466363
$ws.Range("F27").Value = 14.4903184429
$ws.Range("D28").Value = "Brassicasterol"
$ws.Range("E28").Value = 0.00000159663061052
$ws.Range("F28").Value = 6.68035574293
$ws.Range("D29").Value = "Desmosterol"
$ws.Range("E29").Value = 0.00000993571579463
$ws.Range("F29").Value = 26.6135450932
$ws.Range("D30").Value = "Cholesterol"
$ws.Range("E30").Value = 0.0149154483372
$ws.Range("F30").Value = 4.59624688049
$ws.Range("D31").Value = "Cholestanol"
$ws.Range("E31").Value = 0.00205071001666
$ws.Range("F31").Value = 10.0253132674
$ws.Range("D32").Value = "Dehydrocholesterol"
$ws.Range("E32").Value = 0.00285513138424
$ws.Range("F32").Value = 6.22549391704
$ws.Range("D33").Value = "Ergosterol"
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = "inf"
$ws.Range("D34").Value = "Total"
$ws.Range("E34").Value = 7.72670281262
$ws.Range("F34").Value = 0.110282326655
$ws.Range("D35").Value = "fecales"
$ws.Range("E35").Value = 5.60435010926
$ws.Range("F35").Value = 0.121520056332
$ws.Range("D36").Value = "fitosteroles"
$ws.Range("E36").Value = 0.449994029176
$ws.Range("F36").Value = 0.142093144529
$ws.Range("D37").Value = 0.196462086872
$ws.Range("D38").Value = "Retention"
$ws.Range("E38").Value = "in"
$ws.Range("F38").Value = "sediment"
$ws.Range("G38").Value = "(%)"
$ws.Range("D39").Value = "N"
$ws.Range("D40").Value = "Coprostanol"
$ws.Range("E40").Value = 0.000000367290511908
$ws.Range("F40").Value = 2.22312576685
$ws.Range("D41").Value = "Epicoprostanol"
$ws.Range("E41").Value = 0.000000434047450237
$ws.Range("F41").Value = 5.9492778485
$ws.Range("D42").Value = "Ethylcoprostanol"
$ws.Range("E42").Value = 0.00000052356514698
$ws.Range("F42").Value = 2.54064387856
$ws.Range("D43").Value = "Coprostanone"
$ws.Range("E43").Value = 0.000000616163588502
$ws.Range("F43").Value = 3.73443419439
$ws.Range("D44").Value = "Coprostane"
$ws.Range("E44").Value = 0.000000000130853609689
$ws.Range("F44").Value = 44.8527572882
$ws.Range("D45").Value = "b-Sitosterol"
$ws.Range("E45").Value = 0.00000857022958269
$ws.Range("F45").Value = 2.8797749267
$ws.Range("D46").Value = "γ-Sitosterol"
$ws.Range("E46").Value = 0.000000115987564817
$ws.Range("F46").Value = 4.39043954356
$ws.Range("D47").Value = "Stigmasterol"
$ws.Range("E47").Value = 0.00000357036444705
$ws.Range("F47").Value = 3.29529338057
$ws.Range("D48").Value = "Stigmastanol"
$ws.Range("E48").Value = 0.0000037693182222
$ws.Range("F48").Value = 3.14841508354
$ws.Range("D49").Value = "Campesterol"
$ws.Range("E49").Value = 0.00000386238112646
$ws.Range("F49").Value = 3.40394931002
$ws.Range("D50").Value = "Campestanol"
$ws.Range("E50").Value = 0.000000165523052131
$ws.Range("F50").Value = 2.57054458367
$ws.Range("D51").Value = "Brassicasterol"
$ws.Range("E51").Value = 0.0000000515977708613
$ws.Range("F51").Value = 6.15644176132
$ws.Range("D52").Value = "Desmosterol"
$ws.Range("E52").Value = 0.0000000115005432031
$ws.Range("F52").Value = 16.7360107268
$ws.Range("D53").Value = "Cholesterol"
$ws.Range("E53").Value = 0.0000135465653221
$ws.Range("F53").Value = 1.56745822565
$ws.Range("D54").Value = "Cholestanol"
$ws.Range("E54").Value = 0.000000268113466969
$ws.Range("F54").Value = 6.09006213808
$ws.Range("D55").Value = "Dehydrocholesterol"
$ws.Range("E55").Value = 0.0000026662987655
$ws.Range("F55").Value = 1.74677976836
$ws.Range("D56").Value = "Ergosterol"
$ws.Range("E56").Value = 0.000000561754761268
$ws.Range("F56").Value = 0.923471867965
$ws.Range("D57").Value = "Total"
$ws.Range("E57").Value = 0.0541453859137
$ws.Range("F57").Value = 0.00201071614638
$ws.Range("D58").Value = "fecales"
$ws.Range("E58").Value = 0.00301717807049
$ws.Range("F58").Value = 0.00259104730115
$ws.Range("D59").Value = "fitosteroles"
$ws.Range("E59").Value = 0.0311895761621
$ws.Range("F59").Value = 0.00224096633539

# Scientific-notation formatting for the very small "N retention" values.
$ws.Range("E21").NumberFormat = "0.00E+000"
$ws.Range("E27").NumberFormat = "0.00E+000"
$ws.Range("E28").NumberFormat = "0.00E+000"
$ws.Range("E29").NumberFormat = "0.00E+000"
$ws.Range("E40").NumberFormat = "0.00E+000"
$ws.Range("E41").NumberFormat = "0.00E+000"
$ws.Range("E42").NumberFormat = "0.00E+000"
$ws.Range("E43").NumberFormat = "0.00E+000"
$ws.Range("E44").NumberFormat = "0.00E+000"
$ws.Range("E45").NumberFormat = "0.00E+000"
$ws.Range("E46").NumberFormat = "0.00E+000"
$ws.Range("E47").NumberFormat = "0.00E+000"
$ws.Range("E48").NumberFormat = "0.00E+000"
$ws.Range("E49").NumberFormat = "0.00E+000"
$ws.Range("E50").NumberFormat = "0.00E+000"
$ws.Range("E51").NumberFormat = "0.00E+000"
$ws.Range("E52").NumberFormat = "0.00E+000"
$ws.Range("E53").NumberFormat = "0.00E+000"
$ws.Range("E54").NumberFormat = "0.00E+000"
$ws.Range("E55").NumberFormat = "0.00E+000"
$ws.Range("E56").NumberFormat = "0.00E+000"

# --- sed sheet: corrected reactivities (x 10/9) + scratch cell ---
$ws = $wb.Worksheets.Item("sed")
$ws.Range("E2").Value = 348.998690737612
$ws.Range("F2").Value = 99.5394673723764
$ws.Range("G2").Value = 53.4844365945727
$ws.Range("H2").Value = 40.3582522194336
$ws.Range("I2").Value = 1.51084270391904
$ws.Range("J2").Value = 29.9714649811107
$ws.Range("K2").Value = 0.142917229408645
$ws.Range("L2").Value = 3.23805120623662
$ws.Range("M2").Value = 13.7322942290163
$ws.Range("N2").Value = 3.68612421062472
$ws.Range("O2").Value = 0.07395518
$ws.Range("P2").Value = 0.00851811
$ws.Range("Q2").Value = 0.211174314021521
$ws.Range("R2").Value = 54.7493368202176
$ws.Range("S2").Value = 16.4187998059137
$ws.Range("T2").Value = 14.1951214849675
$ws.Range("U2").Value = 0.19810494785578
$ws.Range("A5").Value = 314.0988

# --- View state: selections per sheet, and which tab is active ---
$wsWFL = $wb.Worksheets.Item("WFL")
$wsWFL.Range("E2").Select()

$wsDW = $wb.Worksheets.Item("DW")
$wsDW.Activate()
$wsDW.Range("F54").Select()

$wsSed = $wb.Worksheets.Item("sed")
$wsSed.Activate()
$wsSed.Range("C14").Select()

$wsSheet3 = $wb.Worksheets.Item("Sheet3")
$wsSheet3.Activate()
$wsSheet3.Range("K25").Select()
